# "Diego no me llores todo esta bien"
# Update the MODCOD table: row 6 (16APSK 8/9) CN_req value corrected from
# 12.08 to 13.08, and the temporary yellow highlight on that row is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Correct the CN_req value for row ID 6 (16APSK 8/9)
$ws.Range("C7").Value = 13.08

# Remove the yellow highlight fill from row 7 (A7:D7) by copying the
# (unfilled) formatting from the row below it, which uses the same borders.
$ws.Range("A8:D8").Copy() | Out-Null
$ws.Range("A7:D7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Restore the selection to match the saved workbook state
$ws.Range("K15").Select() | Out-Null
